{"js": "// The document body contains one title paragraph (\"<date> <weekday>\")\n// followed by a 5x5 table of \"AA\u00d7BB=\" multiplication prompts (rows 2-4\n// of every 5-row block are intentionally blank for student work).\n// `context.document.body.paragraphs` walks the whole body in document\n// order, including every table-cell paragraph, so we can address each\n// populated paragraph positionally exactly as the OOXML diff does and\n// swap its text in place (keeping the existing run formatting/fonts).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map from the paragraph's current (old) text to its replacement, in the\n// same order the diff rewrites each <w:t>.\nconst replacements = [\n  [\"2025-10-03 Friday\", \"2025-10-04 Saturday\"],\n  [\"73\u00d728=\", \"98\u00d796=\"],\n  [\"76\u00d731=\", \"93\u00d771=\"],\n  [\"67\u00d783=\", \"62\u00d743=\"],\n  [\"71\u00d761=\", \"43\u00d744=\"],\n  [\"28\u00d739=\", \"42\u00d759=\"],\n  [\"94\u00d753=\", \"62\u00d767=\"],\n  [\"37\u00d781=\", \"46\u00d711=\"],\n  [\"80\u00d754=\", \"42\u00d751=\"],\n  [\"20\u00d765=\", \"77\u00d798=\"],\n  [\"24\u00d719=\", \"61\u00d781=\"],\n  [\"35\u00d750=\", \"72\u00d791=\"],\n  [\"30\u00d742=\", \"33\u00d725=\"],\n  [\"19\u00d735=\", \"30\u00d787=\"],\n  [\"65\u00d751=\", \"72\u00d788=\"],\n  [\"68\u00d797=\", \"26\u00d717=\"],\n  [\"88\u00d759=\", \"97\u00d795=\"],\n  [\"94\u00d786=\", \"27\u00d737=\"],\n  [\"86\u00d776=\", \"75\u00d760=\"],\n  [\"67\u00d765=\", \"18\u00d739=\"],\n  [\"93\u00d750=\", \"29\u00d775=\"],\n  [\"37\u00d772=\", \"60\u00d783=\"],\n  [\"66\u00d734=\", \"69\u00d793=\"],\n  [\"81\u00d745=\", \"76\u00d777=\"],\n  [\"54\u00d723=\", \"59\u00d783=\"],\n  [\"62\u00d715=\", \"97\u00d737=\"],\n];\n\nconst items = paragraphs.items;\nlet repIdx = 0;\nfor (let i = 0; i < items.length && repIdx < replacements.length; i++) {\n  const [oldText, newText] = replacements[repIdx];\n  if (items[i].text === oldText) {\n    items[i].insertText(newText, \"Replace\");\n    repIdx++;\n  }\n}\n\nawait context.sync();\n\nif (repIdx !== replacements.length) {\n  throw new Error(\n    `Only matched ${repIdx} of ${replacements.length} expected paragraphs`\n  );\n}\n", "ps1": "# The document body is a title paragraph (\"<date> <weekday>\") followed by\n# a 5x5 table of \"AA\u00d7BB=\" multiplication prompts (rows 2-4 of every 5-row\n# block are intentionally left blank for student work). Every populated\n# cell/title text is unique in the document, so each old value can be\n# located and swapped for its replacement with a simple Find/Replace pass\n# over the whole document content range.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-10-03 Friday\", \"2025-10-04 Saturday\"),\n    @(\"73\u00d728=\", \"98\u00d796=\"),\n    @(\"76\u00d731=\", \"93\u00d771=\"),\n    @(\"67\u00d783=\", \"62\u00d743=\"),\n    @(\"71\u00d761=\", \"43\u00d744=\"),\n    @(\"28\u00d739=\", \"42\u00d759=\"),\n    @(\"94\u00d753=\", \"62\u00d767=\"),\n    @(\"37\u00d781=\", \"46\u00d711=\"),\n    @(\"80\u00d754=\", \"42\u00d751=\"),\n    @(\"20\u00d765=\", \"77\u00d798=\"),\n    @(\"24\u00d719=\", \"61\u00d781=\"),\n    @(\"35\u00d750=\", \"72\u00d791=\"),\n    @(\"30\u00d742=\", \"33\u00d725=\"),\n    @(\"19\u00d735=\", \"30\u00d787=\"),\n    @(\"65\u00d751=\", \"72\u00d788=\"),\n    @(\"68\u00d797=\", \"26\u00d717=\"),\n    @(\"88\u00d759=\", \"97\u00d795=\"),\n    @(\"94\u00d786=\", \"27\u00d737=\"),\n    @(\"86\u00d776=\", \"75\u00d760=\"),\n    @(\"67\u00d765=\", \"18\u00d739=\"),\n    @(\"93\u00d750=\", \"29\u00d775=\"),\n    @(\"37\u00d772=\", \"60\u00d783=\"),\n    @(\"66\u00d734=\", \"69\u00d793=\"),\n    @(\"81\u00d745=\", \"76\u00d777=\"),\n    @(\"54\u00d723=\", \"59\u00d783=\"),\n    @(\"62\u00d715=\", \"97\u00d737=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
